{"js": "// Update the group-count specification in the \"Trabalho1-PC2\" document:\n//   \"Como s\u00e3o 37 alunos, apenas 1 grupo ter\u00e1 5 pessoas. \"\n// becomes\n//   \"Como s\u00e3o 38 alunos, ser\u00e3o apenas 2 grupos de 3 pessoas. \"\n//\n// The first part (\"Como s\u00e3o 37 alunos,\") is not bold; the second part\n// (\" apenas 1 grupo ter\u00e1 5 pessoas. \") is bold. We preserve that split\n// and its formatting while updating the wording/numbers.\n\nconst body = context.document.body;\n\n// 1) Fix the student count: 37 -> 38 (keeps the surrounding non-bold run).\nconst countMatches = body.search(\"37\", { matchCase: true, matchWholeWord: false });\ncountMatches.load(\"items\");\nawait context.sync();\n\nif (countMatches.items.length > 0) {\n  countMatches.items[0].insertText(\"38\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Replace the bold sentence about the number of groups/people.\nconst groupMatches = body.search(\" apenas 1 grupo ter\u00e1 5 pessoas. \", { matchCase: true });\ngroupMatches.load(\"items\");\nawait context.sync();\n\nif (groupMatches.items.length > 0) {\n  groupMatches.items[0].insertText(\n    \" ser\u00e3o apenas 2 grupos de 3 pessoas. \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Update the group-count specification in the \"Trabalho1-PC2\" document:\n#   \"Como s\u00e3o 37 alunos, apenas 1 grupo ter\u00e1 5 pessoas. \"\n# becomes\n#   \"Como s\u00e3o 38 alunos, ser\u00e3o apenas 2 grupos de 3 pessoas. \"\n#\n# The first part (\"Como s\u00e3o 37 alunos,\") is not bold; the second part\n# (\" apenas 1 grupo ter\u00e1 5 pessoas. \") is bold. We preserve that split\n# and its formatting while updating the wording/numbers.\n\n$d = $word.ActiveDocument\n\n# 1) Fix the student count: 37 -> 38 (keeps the surrounding non-bold run).\n$countRange = $d.Content\n$countRange.Find.Text = \"37\"\n$countFound = $countRange.Find.Execute()\nif ($countFound) {\n    $countRange.Text = \"38\"\n}\n\n# 2) Replace the bold sentence about the number of groups/people.\n$groupRange = $d.Content\n$groupRange.Find.Text = \" apenas 1 grupo ter\u00e1 5 pessoas. \"\n$groupFound = $groupRange.Find.Execute()\nif ($groupFound) {\n    $groupRange.Text = \" ser\u00e3o apenas 2 grupos de 3 pessoas. \"\n}\n"}
